$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (dates) on the existing rows carries a bold/bordered/centered date-format style.
# Copy that formatting from A51 onto the new rows before writing the date values into them.
$ws.Range("A51").Copy()
$ws.Range("A52:A54").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Append new data rows 52-54 (columns A:Y) right after the existing A2:Y51 block,
# in the same order as the source rows (oldest to newest).
$newRowNumbers = @(52, 53, 54)
$newRowValues = @(
    @(45482, 619.581582, 212.6303177275, 0, 0.04040754032, 0, 103.33119885, 0, 207.3473014638, 0, 0.0513980454492, 0, 0, 139.83409027648, 52.504308753, 0, 0.0000021792, 0, 0, 0, 300.5052635923029, 0, 0, 0, 0),
    @(45483, 616.121851254, 215.0154881675, 0, 0.0417564316, 0, 103.15141425, 0, 208.1740164058, 0, 0.0541986074172, 0, 0, 146.7353414976, 53.296596639, 0, 0.0000022104, 0, 0, 0, 319.82254533025, 0, 0, 0, 0),
    @(45484, 612.0024075436, 214.9128703695, 0, 13.14671106129, 0, 107.0617293, 0, 199.0155417369, 0, 20.9920014050975, 0, 0, 140.97313174016, 53.30675417599999, 0, 0.0000021336, 0, 0, 0, 297.0511801027362, 0, 0, 0, 0)
)

for ($i = 0; $i -lt $newRowNumbers.Length; $i++) {
    $r = $newRowNumbers[$i]
    $rowVals = $newRowValues[$i]
    for ($col = 1; $col -le $rowVals.Length; $col++) {
        $ws.Cells.Item($r, $col).Value = $rowVals[$col - 1]
    }
}
